$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the style used by the
# existing header cells (e.g. H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the I and J columns (rows 2-48) with the new data values.
$iValues = @(2,2,7,5,4,4,2,7,5,7,6,7,9,7,6,8,4,7,7,4,11,7,3,6,9,2,8,6,8,6,7,5,4,8,5,4,8,5,3,10,4,2,7,7,7,1,1)
$jValues = @(5,6,8,7,7,6,6,8,6,8,6,8,9,7,7,8,6,8,8,6,11,8,6,7,9,4,9,8,9,8,8,7,4,8,7,7,9,6,7,10,8,6,9,8,8,3,2)

for ($r = 2; $r -le 48; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
